# Fruta / hortaliza, semanal
# Insert a new weekly record at the top of the data block (row 220),
# shifting all subsequent rows down by one. The new row carries the
# latest observation for "Vega Modelo de Temuco - Ciboulette".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row, pushing old rows 220..332 down to 221..333.
$ws.Rows("220:220").Insert()

# Populate the newly-inserted row 220 with the new weekly observation.
$ws.Cells.Item(220, 1).Value = 10                                  # A: Mercado ID
$ws.Cells.Item(220, 2).Value = "Vega Modelo de Temuco"              # B: Mercado
$ws.Cells.Item(220, 3).Value = "La Araucanía"                       # C: Región
$ws.Cells.Item(220, 4).Value = 44960                                # D: Fecha
$ws.Cells.Item(220, 5).Value = 9                                    # E: Codreg
$ws.Cells.Item(220, 6).Value = 100112039                            # F: Categoría ID
$ws.Cells.Item(220, 7).Value = "Ciboulette"                         # G: Categoría
$ws.Cells.Item(220, 8).Value = "Sin especificar"                    # H: Variedad
$ws.Cells.Item(220, 9).Value = "Primera"                            # I: Calidad
$ws.Cells.Item(220, 10).Value = 100                                 # J: Volumen
$ws.Cells.Item(220, 11).Value = 5000                                # K: Precio mínimo
$ws.Cells.Item(220, 12).Value = 6000                                # L: Precio máximo
$ws.Cells.Item(220, 13).Value = 5500                                # M: Precio promedio ponderado
$ws.Cells.Item(220, 14).Value = "$/docena de atados"                # N: Unidad de comercialización
$ws.Cells.Item(220, 15).Value = "Provincia de Cautín"                # O: Origen
$ws.Cells.Item(220, 16).Value = 1833                                # P: Precio $/Kg
$ws.Cells.Item(220, 17).Value = 3                                   # Q: Kg o Unidades
$ws.Cells.Item(220, 18).Value = "Hortaliza"                         # R: Clasificación

# Keep the date format consistent with the rest of column D.
$ws.Cells.Item(220, 4).NumberFormat = $ws.Cells.Item(221, 4).NumberFormat
